$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 8).Value = 0.1067320294967284

$ws.Cells.Item(3, 2).Value = 0.001931220187241999
$ws.Cells.Item(3, 3).Value = 0.0007391134248619421
$ws.Cells.Item(3, 4).Value = 4.533472944071591
$ws.Cells.Item(3, 5).Value = 0.1215017266797351
$ws.Cells.Item(3, 6).Value = 0.0004825777530472764
$ws.Cells.Item(3, 7).Value = 0.003379862621436722
$ws.Cells.Item(3, 8).Value = 0.1086632496839704

$ws.Cells.Item(4, 2).Value = 0.01473899882116804
$ws.Cells.Item(4, 3).Value = 0.001038360306544556
$ws.Cells.Item(4, 4).Value = 13.4161897427049
$ws.Cells.Item(4, 5).Value = 0.0450062568060761
$ws.Cells.Item(4, 6).Value = 0.01270384119112573
$ws.Cells.Item(4, 7).Value = 0.01677415645121035
$ws.Cells.Item(4, 8).Value = 0.1214710283178964

$ws.Cells.Item(5, 2).Value = 0.02371281270418166
$ws.Cells.Item(5, 3).Value = 0.002015050717453312
$ws.Cells.Item(5, 4).Value = 15.86525303363705
$ws.Cells.Item(5, 5).Value = 0.02119810384399553
$ws.Cells.Item(5, 6).Value = 0.01976337128720295
$ws.Cells.Item(5, 7).Value = 0.02766225412116036
$ws.Cells.Item(5, 8).Value = 0.1304448422009101

$ws.Cells.Item(6, 2).Value = 0.05352717466671876
$ws.Cells.Item(6, 3).Value = 0.00882930559469961
$ws.Cells.Item(6, 4).Value = 11.71240800195856
$ws.Cells.Item(6, 5).Value = 0.2235669819378172
$ws.Cells.Item(6, 6).Value = 0.03622199909801526
$ws.Cells.Item(6, 7).Value = 0.07083235023542224
$ws.Cells.Item(6, 8).Value = 0.1602592041634472

$ws.Cells.Item(7, 2).Value = 0.03054790605545981
$ws.Cells.Item(7, 3).Value = 0.007279738129611572
$ws.Cells.Item(7, 4).Value = 7.743486950189914
$ws.Cells.Item(7, 5).Value = 0.09686849510714998
$ws.Cells.Item(7, 6).Value = 0.01627984332713598
$ws.Cells.Item(7, 7).Value = 0.04481596878378363
$ws.Cells.Item(7, 8).Value = 0.1372799355521882

$ws.Cells.Item(8, 2).Value = 0.06546333193498356
$ws.Cells.Item(8, 3).Value = 0.004290622316217816
$ws.Cells.Item(8, 4).Value = 9.519687227143299
$ws.Cells.Item(8, 5).Value = 0.0909594765957963
$ws.Cells.Item(8, 6).Value = 0.05705384216555942
$ws.Cells.Item(8, 7).Value = 0.07387282170440769
$ws.Cells.Item(8, 8).Value = 0.172195361431712

$ws.Cells.Item(9, 2).Value = 0.05587805816401457
$ws.Cells.Item(9, 3).Value = 0.004192425893160631
$ws.Cells.Item(9, 4).Value = 9.58108517232893
$ws.Cells.Item(9, 5).Value = 0.1135099855942992
$ws.Cells.Item(9, 6).Value = 0.04766102976985907
$ws.Cells.Item(9, 7).Value = 0.06409508655817008
$ws.Cells.Item(9, 8).Value = 0.162610087660743

$ws.Cells.Item(10, 2).Value = -0.1067320294967284
$ws.Cells.Item(10, 3).Value = 0.000620563534472687
$ws.Cells.Item(10, 4).Value = -219.7236669231151
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = -0.1079483173860676
$ws.Cells.Item(10, 7).Value = -0.1055157416073892

$ws.Cells.Item(11, 2).Value = -0.05764776983284882
$ws.Cells.Item(11, 3).Value = 0.0006668011590463343
$ws.Cells.Item(11, 4).Value = -101.7770820903605
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = -0.05895468212092519
$ws.Cells.Item(11, 7).Value = -0.05634085754477245
$ws.Cells.Item(11, 8).Value = 0.04908425966387958

$ws.Cells.Item(12, 2).Value = -0.04730157051836866
$ws.Cells.Item(12, 3).Value = 0.0006651746031243923
$ws.Cells.Item(12, 4).Value = -83.14165977590672
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = -0.04860529488269105
$ws.Cells.Item(12, 7).Value = -0.04599784615404628
$ws.Cells.Item(12, 8).Value = 0.05943045897835974

$ws.Cells.Item(13, 2).Value = -0.04478056843773769
$ws.Cells.Item(13, 3).Value = 0.0006560891359623006
$ws.Cells.Item(13, 4).Value = -76.00661974230002
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = -0.04606648549955093
$ws.Cells.Item(13, 7).Value = -0.04349465137592445
$ws.Cells.Item(13, 8).Value = 0.06195146105899071

$ws.Cells.Item(14, 2).Value = -0.03840347378077935
$ws.Cells.Item(14, 3).Value = 0.0006570525705733564
$ws.Cells.Item(14, 4).Value = -64.03748314297083
$ws.Cells.Item(14, 5).Value = (1.783276166231692 * [Math]::Pow(10, -252))
$ws.Cells.Item(14, 6).Value = -0.0396912791738007
$ws.Cells.Item(14, 7).Value = -0.03711566838775799
$ws.Cells.Item(14, 8).Value = 0.06832855571594905

$ws.Cells.Item(15, 2).Value = -0.03450124015210326
$ws.Cells.Item(15, 3).Value = 0.0006470633406355668
$ws.Cells.Item(15, 4).Value = -56.80231489329345
$ws.Cells.Item(15, 5).Value = (3.409868282245064 * [Math]::Pow(10, -100))
$ws.Cells.Item(15, 6).Value = -0.03576946688869777
$ws.Cells.Item(15, 7).Value = -0.03323301341550872
$ws.Cells.Item(15, 8).Value = 0.07223078934462515

$ws.Cells.Item(16, 2).Value = -0.03062336409569633
$ws.Cells.Item(16, 3).Value = 0.0006264685092911343
$ws.Cells.Item(16, 4).Value = -50.85524527484909
$ws.Cells.Item(16, 5).Value = (2.126299469271832 * [Math]::Pow(10, -28))
$ws.Cells.Item(16, 6).Value = -0.03185122559146743
$ws.Cells.Item(16, 7).Value = -0.02939550259992525
$ws.Cells.Item(16, 8).Value = 0.07610866540103206

$ws.Cells.Item(17, 2).Value = -0.02697592631395112
$ws.Cells.Item(17, 3).Value = 0.0006364843789788885
$ws.Cells.Item(17, 4).Value = -42.63446552602955
$ws.Cells.Item(17, 5).Value = 0.001428210544352365
$ws.Cells.Item(17, 6).Value = -0.0282234186675458
$ws.Cells.Item(17, 7).Value = -0.02572843396035644
$ws.Cells.Item(17, 8).Value = 0.07975610318277729

$ws.Cells.Item(18, 2).Value = -0.02364257366137205
$ws.Cells.Item(18, 3).Value = 0.0006417453261673837
$ws.Cells.Item(18, 4).Value = -39.10241848913063
$ws.Cells.Item(18, 5).Value = (6.141939363711622 * [Math]::Pow(10, -6))
$ws.Cells.Item(18, 6).Value = -0.02490037736511034
$ws.Cells.Item(18, 7).Value = -0.02238476995763376
$ws.Cells.Item(18, 8).Value = 0.08308945583535635

$ws.Cells.Item(19, 2).Value = -0.0212238803111161
$ws.Cells.Item(19, 3).Value = 0.0006565116750686046
$ws.Cells.Item(19, 4).Value = -35.51283786107579
$ws.Cells.Item(19, 5).Value = (5.767794762741496 * [Math]::Pow(10, -7))
$ws.Cells.Item(19, 6).Value = -0.02251062585012075
$ws.Cells.Item(19, 7).Value = -0.01993713477211145
$ws.Cells.Item(19, 8).Value = 0.08550814918561231

$ws.Cells.Item(20, 2).Value = -0.01789779183669642
$ws.Cells.Item(20, 3).Value = 0.0006686760352865118
$ws.Cells.Item(20, 4).Value = -29.67776527464417
$ws.Cells.Item(20, 5).Value = 0.00644394930116758
$ws.Cells.Item(20, 6).Value = -0.0192083792674242
$ws.Cells.Item(20, 7).Value = -0.01658720440596864
$ws.Cells.Item(20, 8).Value = 0.08883423766003198

$ws.Cells.Item(21, 2).Value = -0.01329161328694966
$ws.Cells.Item(21, 3).Value = 0.0006935252583323473
$ws.Cells.Item(21, 4).Value = -21.99787837730142
$ws.Cells.Item(21, 5).Value = 0.04180328025044439
$ws.Cells.Item(21, 6).Value = -0.01465090458715945
$ws.Cells.Item(21, 7).Value = -0.01193232198673987
$ws.Cells.Item(21, 8).Value = 0.09344041620977873

$ws.Cells.Item(22, 2).Value = -0.009893825314779151
$ws.Cells.Item(22, 3).Value = 0.0006773195974083565
$ws.Cells.Item(22, 4).Value = -16.49430368660508
$ws.Cells.Item(22, 5).Value = 0.01684991611631484
$ws.Cells.Item(22, 6).Value = -0.01122135387004944
$ws.Cells.Item(22, 7).Value = -0.008566296759508851
$ws.Cells.Item(22, 8).Value = 0.09683820418194924

$ws.Cells.Item(23, 2).Value = -0.005748127082400781
$ws.Cells.Item(23, 3).Value = 0.0006625870551790332
$ws.Cells.Item(23, 4).Value = -9.597333220678442
$ws.Cells.Item(23, 5).Value = 0.03970129649341818
$ws.Cells.Item(23, 6).Value = -0.007046780153325155
$ws.Cells.Item(23, 7).Value = -0.004449474011476409
$ws.Cells.Item(23, 8).Value = 0.1009839024143276

$ws.Cells.Item(24, 2).Value = -0.001734901941562419
$ws.Cells.Item(24, 3).Value = 0.0006444367190912092
$ws.Cells.Item(24, 4).Value = -4.604185688435824
$ws.Cells.Item(24, 5).Value = 0.002862226354514899
$ws.Cells.Item(24, 6).Value = -0.002997980773010808
$ws.Cells.Item(24, 7).Value = -0.0004718231101140314
$ws.Cells.Item(24, 8).Value = 0.104997127555166

$ws.Cells.Item(25, 2).Value = (2.396341798272552 * [Math]::Pow(10, -5))
$ws.Cells.Item(25, 3).Value = 0.0006204025367496545
$ws.Cells.Item(25, 4).Value = -2.57339756687702
$ws.Cells.Item(25, 5).Value = 0.1187098307705495
$ws.Cells.Item(25, 6).Value = -0.001192009003564632
$ws.Cells.Item(25, 7).Value = 0.001239935839530083
$ws.Cells.Item(25, 8).Value = 0.1067559929147111

$ws.Cells.Item(26, 2).Value = 0.068371445100081
$ws.Cells.Item(26, 3).Value = 0.001013548164169905
$ws.Cells.Item(26, 4).Value = 66.55032712428445
$ws.Cells.Item(26, 5).Value = 0.1000347722897597
$ws.Cells.Item(26, 6).Value = 0.06638491861860862
$ws.Cells.Item(26, 7).Value = 0.07035797158155337
$ws.Cells.Item(26, 8).Value = 0.1751034745968094
